# Saldo.xlsx update:
# Three accounts on the "Export" sheet receive new (lower) Saldo balances.
# The sheet is kept sorted by Saldo descending, so after the value updates
# the whole data range is re-sorted to restore that order.
#
#   Conta       Nome       old Saldo   -> new Saldo
#   008054713   MODULAR    90000       -> 5000
#   000806386   FERNANDA   70186.25    -> 186.25
#   004329030   DANIELA    50000       -> 940.23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last used row of the data block (header in row 1, data starts
# row 2). The sheet also carries a trailing blank row plus a footer/notes
# row below the data (text in column A only, nothing in column C), so
# anchor the "last row" search on column C (the numeric Saldo column) to
# land exactly on the last real data row instead of the notes row.
$lastDataRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

# --- 1. Update the three balances in place (by account number) ---
for ($r = 2; $r -le $lastDataRow; $r++) {
    $conta = $ws.Cells.Item($r, 1).Value()
    if ($conta -eq "008054713") {
        $ws.Cells.Item($r, 3).Value = 5000
    }
    elseif ($conta -eq "000806386") {
        $ws.Cells.Item($r, 3).Value = 186.25
    }
    elseif ($conta -eq "004329030") {
        $ws.Cells.Item($r, 3).Value = 940.23
    }
}

# --- 2. Re-sort the data rows (A2:C<lastDataRow>) by Saldo, descending ---
$dataRange = $ws.Range("A2:C" + $lastDataRow)
$sortKey = $ws.Range("C2:C" + $lastDataRow)
$dataRange.Sort($sortKey, 2)
